# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Odin_Profits workbook
# (scheduled runner refresh of currentAveragePrice / LevePrice / LeveProfit columns).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 387.92307
$ws.Range("J38").Value = 1727
$ws.Range("L38").Value = 5181
$ws.Range("N38").Value = -5925
# Row 51
$ws.Range("H51").Value = 16710.375
$ws.Range("J51").Value = 3232
$ws.Range("L51").Value = 3232
$ws.Range("N51").Value = -4200
# Row 69
$ws.Range("H69").Value = 9509.6
$ws.Range("I69").Value = 6779.8
$ws.Range("J69").Value = 12239.4
$ws.Range("K69").Value = 20339.4
$ws.Range("L69").Value = 36718.2
$ws.Range("M69").Value = -19465.4
$ws.Range("N69").Value = -38466.2
# Row 72
$ws.Range("H72").Value = 9509.6
$ws.Range("I72").Value = 6779.8
$ws.Range("J72").Value = 12239.4
$ws.Range("K72").Value = 61018.2
$ws.Range("L72").Value = 110154.6
$ws.Range("M72").Value = -56650.2
$ws.Range("N72").Value = -118890.6
# Row 82
$ws.Range("H82").Value = 1377.1666
$ws.Range("I82").Value = 1443.2
$ws.Range("J82").Value = 1047
$ws.Range("K82").Value = 4329.6
$ws.Range("L82").Value = 3141
$ws.Range("M82").Value = -3923.6
$ws.Range("N82").Value = -3953
# Row 85
$ws.Range("H85").Value = 1377.1666
$ws.Range("I85").Value = 1443.2
$ws.Range("J85").Value = 1047
$ws.Range("K85").Value = 4329.6
$ws.Range("L85").Value = 3141
$ws.Range("M85").Value = -2925.6
$ws.Range("N85").Value = -5949
# Row 99
$ws.Range("H99").Value = 83339816
$ws.Range("I99").Value = 600.125
$ws.Range("J99").Value = 250018260
$ws.Range("K99").Value = 1800.375
$ws.Range("L99").Value = 750054780
$ws.Range("M99").Value = -302.375
$ws.Range("N99").Value = -750057776
# Row 138
$ws.Range("H138").Value = 8513.857
$ws.Range("J138").Value = 8611.191999999999
$ws.Range("L138").Value = 25833.576
$ws.Range("N138").Value = -36113.576

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1979.4509
$ws.Range("I32").Value = 624.5238000000001
$ws.Range("K32").Value = 624.5238000000001
$ws.Range("M32").Value = -337.5238000000001
# Row 39
$ws.Range("H39").Value = 1983.3334
$ws.Range("I39").Value = 1975
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 1975
$ws.Range("L39").Value = 2000
$ws.Range("M39").Value = -1455
$ws.Range("N39").Value = -3040

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4591.466
$ws.Range("I31").Value = 1149.1428
$ws.Range("J31").Value = 5408.288
$ws.Range("K31").Value = 1149.1428
$ws.Range("L31").Value = 5408.288
$ws.Range("M31").Value = -854.1428000000001
$ws.Range("N31").Value = -5998.288
# Row 34
$ws.Range("H34").Value = 4591.466
$ws.Range("I34").Value = 1149.1428
$ws.Range("J34").Value = 5408.288
$ws.Range("K34").Value = 1149.1428
$ws.Range("L34").Value = 5408.288
$ws.Range("M34").Value = -947.1428000000001
$ws.Range("N34").Value = -5812.288

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Range("H59").Value = 3423.5
$ws.Range("I59").Value = 1987
$ws.Range("K59").Value = 5961
$ws.Range("M59").Value = -5421
# Row 63
$ws.Range("H63").Value = 18072
$ws.Range("I63").Value = 7776
$ws.Range("J63").Value = 19542.857
$ws.Range("K63").Value = 23328
$ws.Range("L63").Value = 58628.571
$ws.Range("M63").Value = -22579
$ws.Range("N63").Value = -60126.571
# Row 66
$ws.Range("H66").Value = 18072
$ws.Range("I66").Value = 7776
$ws.Range("J66").Value = 19542.857
$ws.Range("K66").Value = 69984
$ws.Range("L66").Value = 175885.713
$ws.Range("M66").Value = -66240
$ws.Range("N66").Value = -183373.713
# Row 68
$ws.Range("H68").Value = 221647.56
$ws.Range("J68").Value = 317725
$ws.Range("L68").Value = 953175
$ws.Range("N68").Value = -954797
# Row 71
$ws.Range("H71").Value = 221647.56
$ws.Range("J71").Value = 317725
$ws.Range("L71").Value = 2859525
$ws.Range("N71").Value = -2867637
# Row 74
$ws.Range("H74").Value = 29166
$ws.Range("J74").Value = 29166
$ws.Range("L74").Value = 87498
$ws.Range("N74").Value = -89620
# Row 77
$ws.Range("H77").Value = 29166
$ws.Range("J77").Value = 29166
$ws.Range("L77").Value = 262494
$ws.Range("N77").Value = -273102
# Row 93
$ws.Range("H93").Value = 14633
$ws.Range("I93").Value = 20899
$ws.Range("J93").Value = 11500
$ws.Range("K93").Value = 62697
$ws.Range("L93").Value = 34500
$ws.Range("M93").Value = -60825
$ws.Range("N93").Value = -38244
# Row 94
$ws.Range("H94").Value = 5874.75
$ws.Range("I94").Value = 2500
$ws.Range("J94").Value = 6999.6665
$ws.Range("K94").Value = 7500
$ws.Range("L94").Value = 20998.9995
$ws.Range("M94").Value = -6824
$ws.Range("N94").Value = -22350.9995
# Row 98
$ws.Range("H98").Value = 4699
$ws.Range("I98").Value = 8800
$ws.Range("J98").Value = 3878.8
$ws.Range("K98").Value = 26400
$ws.Range("L98").Value = 11636.4
$ws.Range("M98").Value = -24902
$ws.Range("N98").Value = -14632.4
# Row 99
$ws.Range("H99").Value = 7782.222
$ws.Range("I99").Value = 7508
$ws.Range("J99").Value = 8125
$ws.Range("K99").Value = 22524
$ws.Range("L99").Value = 24375
$ws.Range("M99").Value = -20278
$ws.Range("N99").Value = -28867
# Row 101
$ws.Range("H101").Value = 38333
$ws.Range("J101").Value = 38333
$ws.Range("L101").Value = 114999
$ws.Range("N101").Value = -119867
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
# Row 105
$ws.Range("H105").Value = 28132.834
$ws.Range("J105").Value = 28132.834
$ws.Range("L105").Value = 84398.50199999999
$ws.Range("N105").Value = -89640.50199999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 136
$ws.Range("H136").Value = 30746.857
$ws.Range("J136").Value = 30746.857
$ws.Range("L136").Value = 92240.571
$ws.Range("N136").Value = -97340.571

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1002
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1002
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 2004
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -4126
# Row 84
$ws.Range("H84").Value = 1002
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1002
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -20628
